$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.095.75"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "3.483.54"
$ws.Range("E3").Value = "  -5.37%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'607.57"
$ws.Range("E5").Value = "  -6.59%  "
$ws.Range("D6").Value = "'148.62"
$ws.Range("E6").Value = "  -7.95%  "
$ws.Range("D7").Value = "3.481.24"
$ws.Range("E7").Value = "  -5.38%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.63%  "
$ws.Range("E10").Value = "  -5.69%  "
$ws.Range("D11").Value = "'6.88"
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("E12").Value = "  -4.91%  "
$ws.Range("D13").Value = "'0.0000213"
$ws.Range("E13").Value = "  -8.03%  "
$ws.Range("D14").Value = "4.072.27"
$ws.Range("E14").Value = "  -5.35%  "
$ws.Range("D15").Value = "'31.32"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").Value = "3.477.19"
$ws.Range("E16").Value = "  -4.90%  "
$ws.Range("D17").Value = "66.977.47"
$ws.Range("E17").Value = "  -4.06%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "'15.03"
$ws.Range("E20").Value = "  -5.69%  "
$ws.Range("D21").Value = "'444.02"
$ws.Range("E21").Value = "  -5.77%  "
$ws.Range("D22").Value = "'9.03"
$ws.Range("E22").Value = "  -12.68%  "
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").Value = "'77.12"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "3.622.07"
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "'10.14"
$ws.Range("E28").Value = "  -8.30%  "
$ws.Range("D29").Value = "'8.31"
$ws.Range("E29").Value = "  -6.00%  "
$ws.Range("D30").Value = "'2.53"
$ws.Range("E30").Value = "  -4.63%  "
$ws.Range("E31").Value = "  -7.72%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'0.165"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'25.60"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").Value = "'6.12"
$ws.Range("E35").Value = "  -6.18%  "
$ws.Range("E36").Value = "  -7.32%  "
$ws.Range("D37").Value = "3.464.59"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("D38").Value = "'7.97"
$ws.Range("E38").Value = "  -5.13%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "'2.20"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "'170.43"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").Value = "'0.0868"
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("E44").Value = "  -7.65%  "
$ws.Range("D45").Value = "'0.881"
$ws.Range("E45").Value = "  -5.28%  "
$ws.Range("D46").Value = "'45.68"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "'1.24"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").Value = "'2.51"
$ws.Range("E48").Value = "  -11.12%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'26.00"
$ws.Range("E49").Value = "  -10.18%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.55"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("E51").Value = "  -4.32%  "
